$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.448.82"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "'1.924.58"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "'0.740"
$ws.Range("E5").Value = "  +13.74%  "

$ws.Range("D6").Value = "'254.00"
$ws.Range("E6").Value = "  +4.14%  "

$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").Value = "'40.75"
$ws.Range("E8").Value = "  -2.02%  "

$ws.Range("D9").Value = "'0.358"
$ws.Range("E9").Value = "  +4.48%  "

$ws.Range("D10").Value = "'52.65"
$ws.Range("E10").Value = "  +5.00%  "

$ws.Range("D11").Value = "'0.0741"
$ws.Range("E11").Value = "  +4.31%  "

$ws.Range("D12").Value = "'0.0999"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").Value = "'2.200.11"
$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").Value = "'12.71"
$ws.Range("E14").Value = "  +4.99%  "

$ws.Range("D15").Value = "'0.718"
$ws.Range("E15").Value = "  +3.60%  "

$ws.Range("D16").Value = "'1.954.62"
$ws.Range("E16").Value = "  +2.79%  "

$ws.Range("D17").Value = "'4.91"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "'35.437.03"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").Value = "'73.60"
$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("D20").Value = "'0.0₃0840"
$ws.Range("E20").Value = "  +2.92%  "

$ws.Range("D21").Value = "'13.06"
$ws.Range("E21").Value = "  +4.37%  "

$ws.Range("D22").Value = "'242.34"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("E23").Value = "  +7.66%  "

$ws.Range("E24").Value = "  -0.27%  "

$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +2.08%  "

$ws.Range("D26").Value = "'2.34"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").Value = "'167.99"
$ws.Range("E27").Value = "  -1.32%  "

$ws.Range("D28").Value = "'8.73"
$ws.Range("E28").Value = "  +4.57%  "

$ws.Range("E29").Value = "  +8.00%  "

$ws.Range("D30").Value = "'18.89"
$ws.Range("E30").Value = "  +3.74%  "

$ws.Range("D31").Value = "'4.127.85"
$ws.Range("E31").Value = "  +19.43%  "

$ws.Range("B32").Value = "'Filecoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.36"
$ws.Range("E32").Value = "  +5.67%  "

$ws.Range("B33").Value = "'TrustWalletToken"
$ws.Range("C33").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D33").Value = "'1.66"
$ws.Range("E33").Value = "  +24.35%  "

$ws.Range("E34").Value = "  +14.60%  "

$ws.Range("E35").Value = "  +3.38%  "

$ws.Range("E36").Value = "  +3.82%  "

$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("D38").Value = "'0.914"
$ws.Range("E38").Value = "  -2.04%  "

$ws.Range("D39").Value = "'2.05"
$ws.Range("E39").Value = "  +0.56%  "

$ws.Range("B40").Value = "'Aave"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'99.18"
$ws.Range("E40").Value = "  +10.81%  "

$ws.Range("B41").Value = "'InjectiveProtocol"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.37"
$ws.Range("E41").Value = "  +9.49%  "

$ws.Range("E42").Value = "  +4.67%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("E44").Value = "  +3.89%  "

$ws.Range("D45").Value = "'2.50"
$ws.Range("E45").Value = "  +6.26%  "

$ws.Range("D46").Value = "'1.349.77"
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").Value = "'6.70"
$ws.Range("E49").Value = "  +3.11%  "

$ws.Range("D50").Value = "'45.26"
$ws.Range("E50").Value = "  -4.15%  "

$ws.Range("D51").Value = "'2.107.65"
$ws.Range("E51").Value = "  +1.13%  "
